$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "otp"
$ws.Range("B3").Value = 111111

$ws.Range("B3").Select()
